$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing Digikey/Mouser reference for the SD card row (D27)
$ws.Range("D27").Value = "798-DM3CSSF "

# Add the new USB-C connector row
# (write D28 before B28 so the shared-string table gets the two new
# entries in the same order as the reference workbook: 640-USB4085-GF-A
# first, then USB C Buchse)
$ws.Range("A28").Value = 1
$ws.Range("D28").Value = "640-USB4085-GF-A "
$ws.Range("B28").Value = "USB C Buchse"

# Update view state to match the author's saved selection/scroll position
$ws.Range("A29").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
